$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn on AutoFilter for the table range (do this before the new rows are
# appended so the filter range stays anchored to the original table A1:E4)
$ws.Range("A1:E4").AutoFilter()

# Excel records the filter range as a hidden sheet-scoped defined name
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $ws.Name + "!`$A`$1:`$E`$4")
$fdb.Visible = $false

# Highlight the header row with a solid blue fill (0070C0)
$ws.Range("A1:E1").Interior.Color = 12611584

# Add the extra rows of numbers in column A (reading week log numbers)
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# Update the active selection to the header row
$ws.Range("A1:E1").Select()
